# fix lỗi trong report cơ sở. Thêm cột ghi chú trong báo cáo về chi tiêu
$wb = $excel.ActiveWorkbook

# Helper: write a plain-text value into a cell without Excel's "looks like
# a date/number" auto-conversion kicking in (e.g. "08-08-2024" silently
# becoming a date serial). Force the cell to Text format for the write,
# then drop back to the Normal style so no stray per-cell style survives.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet 1: "Đơn sale chính" — currently empty, fill in header + 2 rows
# ---------------------------------------------------------------------
$wsSaleChinh = $wb.Worksheets.Item(1)

$headers1 = @("Tiền tố","Mã dịch vụ","Ngày thực hiện","Cơ sở","Khách hàng","Nguồn khách","Tên dịch vụ","Đơn giá gốc","Sale phụ","Upsale","Đơn giá","Đã thanh toán","Tỉ lệ chiết khấu sale chính","Chiết khấu sale chính")
for ($i = 0; $i -lt $headers1.Count; $i++) {
    Set-TextValue $wsSaleChinh.Cells.Item(1, $i + 1) $headers1[$i]
}

$row2Text = @{1="HD-LUXURY"; 3="08-08-2024"; 4="CẦN THƠ"; 5="Nguyễn Bích Thuỳ"; 6="Khách cũ giới thiệu"; 7="Cắt mí"}
$row2Num  = @{2=635; 8=9000000; 11=9000000; 12=9000000; 13=0.1; 14=900000}
foreach ($col in $row2Text.Keys) { Set-TextValue $wsSaleChinh.Cells.Item(2, $col) $row2Text[$col] }
foreach ($col in $row2Num.Keys)  { $wsSaleChinh.Cells.Item(2, $col).Value = $row2Num[$col] }

Set-TextValue $wsSaleChinh.Cells.Item(3, 1) "Tổng"
$row3Num = @{2=1; 8=9000000; 10=0; 11=9000000; 12=9000000; 13=0; 14=900000}
foreach ($col in $row3Num.Keys) { $wsSaleChinh.Cells.Item(3, $col).Value = $row3Num[$col] }

# ---------------------------------------------------------------------
# Sheet 2: "Đơn 1 bác sĩ" — insert 6 new detail rows before the "Tổng"
# row (which moves from row 7 to row 13), then refresh the totals.
# ---------------------------------------------------------------------
$wsBacSi = $wb.Worksheets.Item(2)

# Move the existing "Tổng" row (currently row 7) down to row 13 first.
for ($k = 0; $k -lt 6; $k++) {
    $wsBacSi.Rows.Item(7).Insert() | Out-Null
}

$newRowsText = @{
    7  = @{1="HD-LUXURY"; 3="08-04-2024"; 4="SÓC TRĂNG"; 5="nguyễn thị mỹ trinh"; 6="Cá nhân"; 7="nhấn đồng tiền"}
    8  = @{1="HD-LUXURY"; 3="08-04-2024"; 4="SÓC TRĂNG"; 5="nguyễn thị lệ trang"; 6="Cá nhân"; 7="Cắt mí"}
    9  = @{1="HD-LUXURY"; 3="08-05-2024"; 4="CẦN THƠ"; 5="Phạm Thị Trúc Lài"; 6="CTV"; 7="Phun mày"; 9="Đỗ Thị Huyền Trân"}
    10 = @{1="HD-LUXURY"; 3="08-06-2024"; 4="CẦN THƠ"; 5="Võ Thị Thuỳ Trang"; 6="CTV"; 7="Phun mày"; 9="Đỗ Thị Huyền Trân"}
    11 = @{1="HD-LUXURY"; 3="08-07-2024"; 4="CẦN THƠ"; 5="Nguyễn Thị Thắm"; 6="Cá nhân"; 7="Tiêm botox"}
    12 = @{1="HD-LUXURY"; 3="08-08-2024"; 4="CẦN THƠ"; 5="Nguyễn Bích Thuỳ"; 6="Khách cũ giới thiệu"; 7="Cắt mí"}
}

$newRowsNum = @{
    7  = @{2=626; 8=7000000;  11=7000000;  12=7000000;  13=0.1;  14=700000}
    8  = @{2=628; 8=4000000;  11=4000000;  12=4000000;  13=0.1;  14=400000}
    9  = @{2=629; 8=500000;   10=500000;   11=1000000;  12=1000000; 13=0.1; 14=100000}
    10 = @{2=631; 8=500000;   10=1000000;  11=1500000;  12=1500000; 13=0.1; 14=150000}
    11 = @{2=632; 8=2000000;  11=2000000;  12=2000000;  13=0.08; 14=160000}
    12 = @{2=635; 8=9000000;  11=9000000;  12=9000000;  13=0.08; 14=720000}
}

foreach ($r in $newRowsText.Keys) {
    $cols = $newRowsText[$r]
    foreach ($col in $cols.Keys) {
        Set-TextValue $wsBacSi.Cells.Item($r, $col) $cols[$col]
    }
}
foreach ($r in $newRowsNum.Keys) {
    $cols = $newRowsNum[$r]
    foreach ($col in $cols.Keys) {
        $wsBacSi.Cells.Item($r, $col).Value = $cols[$col]
    }
}

# Refresh the "Tổng" row, now at row 13.
Set-TextValue $wsBacSi.Cells.Item(13, 1) "Tổng"
$totalRow13 = @{2=11; 8=43100000; 10=7500000; 11=50600000; 12=49600000; 13=0; 14=4578000}
foreach ($col in $totalRow13.Keys) { $wsBacSi.Cells.Item(13, $col).Value = $totalRow13[$col] }

# ---------------------------------------------------------------------
# Sheet 3: "Lương" — refresh totals that ripple from the new rows above.
# ---------------------------------------------------------------------
$wsLuong = $wb.Worksheets.Item(3)

$luongUpdates = @{
    2  = 8
    3  = 280000
    4  = 2285714.285714286
    5  = 900000
    7  = 1778000
    22 = 3
    23 = 857142.8571428573
    27 = 1700000
    32 = 5243714.285714285
    34 = 2557142.857142857
    35 = 9186571.428571429
}

foreach ($r in $luongUpdates.Keys) {
    $wsLuong.Cells.Item($r, 2).Value = $luongUpdates[$r]
}
